$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.221.74'
$ws.Range("E2").Value = '  +0.97%  '
# Row 3
$ws.Range("D3").Value = '2.599.60'
$ws.Range("E3").Value = '  +0.49%  '
# Row 4
$ws.Range("E4").Value = '  -0.04%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.16%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.13%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '
# Row 9
$ws.Range("D9").Value = '2.606.10'
$ws.Range("E9").Value = '  +0.04%  '
# Row 10
$ws.Range("E10").Value = '  -2.26%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.96%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.157'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.66%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.370'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.63%  '
# Row 14
$ws.Range("D14").Value = '3.058.80'
$ws.Range("E14").Value = '  +0.40%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.53%  '
# Row 16
$ws.Range("D16").Value = '60.219.89'
$ws.Range("E16").Value = '  +1.02%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000141'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.55%  '
# Row 18
$ws.Range("D18").Value = '2.605.92'
$ws.Range("E18").Value = '  +0.30%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.45'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.94%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.00%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '346.70'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.64%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.25%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
# Row 24
$ws.Range("E24").Value = '  +4.95%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.89%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.18%  '
# Row 27
$ws.Range("E27").Value = '  -0.25%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.17%  '
# Row 29
$ws.Range("D29").Value = '0.0₃0792'
$ws.Range("E29").Value = '  +2.28%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.68%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.38'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.95%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '164.89'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.97%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.41'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.38%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.57%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.17%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.981'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.54%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.63'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.24%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.70%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '312.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.46%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.66%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.837'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.52%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '134.76'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.80%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0993'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.73%  '
# Row 45
$ws.Range("E45").Value = '  +0.17%  '
# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.84%  '
# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.605'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.80%  '
# Row 48
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.01'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.27%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0550'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.52%  '
# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.00%  '
# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.03%  '
